$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 452.875
$ws.Range("I12").Value = 479.1
$ws.Range("K12").Value = 479.1
$ws.Range("M12").Value = -309.1

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2307
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 2307
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 2307
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -2445

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 6226.3076
$ws.Range("I76").Value = 4298.1665
$ws.Range("K76").Value = 4298.1665
$ws.Range("M76").Value = -3983.1665

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 6226.3076
$ws.Range("I79").Value = 4298.1665
$ws.Range("K79").Value = 4298.1665
$ws.Range("M79").Value = -3206.1665

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H95").Value = 65873.375
$ws.Range("J95").Value = 65873.375
$ws.Range("L95").Value = 65873.375
$ws.Range("N95").Value = -71365.375

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 18657.705
$ws.Range("I135").Value = 1317
$ws.Range("J135").Value = 46669.617
$ws.Range("K135").Value = 11853
$ws.Range("L135").Value = 420026.553
$ws.Range("M135").Value = -9318
$ws.Range("N135").Value = -425096.553

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 8159.5356
$ws.Range("I137").Value = 14982.167
$ws.Range("J137").Value = 3042.5625
$ws.Range("K137").Value = 44946.501
$ws.Range("L137").Value = 9127.6875
$ws.Range("M137").Value = -42396.501
$ws.Range("N137").Value = -14227.6875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 53183
$ws.Range("J24").Value = 53183
$ws.Range("L24").Value = 53183
$ws.Range("N24").Value = -53931

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6220.713
$ws.Range("I32").Value = 5632.1646
$ws.Range("J32").Value = 11779.223
$ws.Range("K32").Value = 5632.1646
$ws.Range("L32").Value = 11779.223
$ws.Range("M32").Value = -5345.1646
$ws.Range("N32").Value = -12353.223

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 5642.8423
$ws.Range("I45").Value = 4293.5835
$ws.Range("J45").Value = 7955.857
$ws.Range("K45").Value = 4293.5835
$ws.Range("L45").Value = 7955.857
$ws.Range("M45").Value = -3916.5835
$ws.Range("N45").Value = -8709.857

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2755.3774
$ws.Range("I61").Value = 2057.0908
$ws.Range("J61").Value = 6169.222
$ws.Range("K61").Value = 2057.0908
$ws.Range("L61").Value = 6169.222
$ws.Range("M61").Value = -1845.0908
$ws.Range("N61").Value = -6593.222

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 183938.88
$ws.Range("I74").Value = 241219.12
$ws.Range("J74").Value = 4938.125
$ws.Range("K74").Value = 241219.12
$ws.Range("L74").Value = 4938.125
$ws.Range("M74").Value = -240345.12
$ws.Range("N74").Value = -6686.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 183938.88
$ws.Range("I77").Value = 241219.12
$ws.Range("J77").Value = 4938.125
$ws.Range("K77").Value = 1206095.6
$ws.Range("L77").Value = 24690.625
$ws.Range("M77").Value = -1201727.6
$ws.Range("N77").Value = -33426.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H100").Value = 53183
$ws.Range("J100").Value = 53183
$ws.Range("L100").Value = 53183
$ws.Range("N100").Value = -55347

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1779.5862
$ws.Range("I102").Value = 1837.3704
$ws.Range("J102").Value = 999.5
$ws.Range("K102").Value = 1837.3704
$ws.Range("L102").Value = 999.5
$ws.Range("M102").Value = -215.3704
$ws.Range("N102").Value = -4243.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2430.5
$ws.Range("I110").Value = 1248.5
$ws.Range("K110").Value = 1248.5
$ws.Range("M110").Value = 796.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2755.3774
$ws.Range("I136").Value = 2057.0908
$ws.Range("J136").Value = 6169.222
$ws.Range("K136").Value = 6171.2724
$ws.Range("L136").Value = 18507.666
$ws.Range("M136").Value = -3621.2724
$ws.Range("N136").Value = -23607.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 27402.37
$ws.Range("I20").Value = 36789
$ws.Range("J20").Value = 1119.8
$ws.Range("K20").Value = 36789
$ws.Range("L20").Value = 1119.8
$ws.Range("M20").Value = -36542
$ws.Range("N20").Value = -1613.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3271.3157
$ws.Range("I99").Value = 2985.0625
$ws.Range("K99").Value = 2985.0625
$ws.Range("M99").Value = -1487.0625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2239.7737
$ws.Range("I134").Value = 1964.3478
$ws.Range("K134").Value = 5893.0434
$ws.Range("M134").Value = -3358.0434

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1935.3077
$ws.Range("I16").Value = 1675.9
$ws.Range("K16").Value = 1675.9
$ws.Range("M16").Value = -1388.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 30899.666
$ws.Range("J50").Value = 30899.666
$ws.Range("L50").Value = 30899.666
$ws.Range("N50").Value = -32149.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 12818.786
$ws.Range("I62").Value = 2996.2
$ws.Range("K62").Value = 2996.2
$ws.Range("M62").Value = -2372.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 12818.786
$ws.Range("I65").Value = 2996.2
$ws.Range("K65").Value = 14981
$ws.Range("M65").Value = -11861

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1935.3077
$ws.Range("I113").Value = 1675.9
$ws.Range("K113").Value = 1675.9
$ws.Range("M113").Value = 494.0999999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 39971.59
$ws.Range("J37").Value = 39971.59
$ws.Range("L37").Value = 119914.77
$ws.Range("N37").Value = -120138.77

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 2842.0435
$ws.Range("I92").Value = 3312.889
$ws.Range("J92").Value = 2539.3572
$ws.Range("K92").Value = 9938.667000000001
$ws.Range("L92").Value = 7618.071599999999
$ws.Range("M92").Value = -8690.667000000001
$ws.Range("N92").Value = -10114.0716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 5706.6
$ws.Range("J137").Value = 7027.6665
$ws.Range("L137").Value = 21082.9995
$ws.Range("N137").Value = -31282.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 207804
$ws.Range("J3").Value = 14000
$ws.Range("L3").Value = 14000
$ws.Range("N3").Value = -14232

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 9362.5
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 9362.5
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 9362.5
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -9640.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 10005
$ws.Range("J12").Value = 10005
$ws.Range("L12").Value = 10005
$ws.Range("N12").Value = -10285

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 24163.334
$ws.Range("J92").Value = 24163.334
$ws.Range("L92").Value = 24163.334
$ws.Range("N92").Value = -27907.334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 874.5
$ws.Range("I97").Value = 874.5
$ws.Range("K97").Value = 874.5
$ws.Range("M97").Value = -378.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 29028
$ws.Range("J98").Value = 29028
$ws.Range("L98").Value = 29028
$ws.Range("N98").Value = -35018

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2282.1843
$ws.Range("I122").Value = 1429.625
$ws.Range("K122").Value = 4288.875
$ws.Range("M122").Value = -1838.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3997.0952
$ws.Range("I132").Value = 3821.6155
$ws.Range("K132").Value = 11464.8465
$ws.Range("M132").Value = -8934.8465

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4559.1
$ws.Range("I46").Value = 6000
$ws.Range("J46").Value = 4399
$ws.Range("K46").Value = 6000
$ws.Range("L46").Value = 4399
$ws.Range("M46").Value = -5812
$ws.Range("N46").Value = -4775

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2237.125
$ws.Range("I61").Value = 2237.125
$ws.Range("K61").Value = 2237.125
$ws.Range("M61").Value = -2035.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 8546.857
$ws.Range("I82").Value = 8513.25
$ws.Range("K82").Value = 8513.25
$ws.Range("M82").Value = -8152.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 8546.857
$ws.Range("I85").Value = 8513.25
$ws.Range("K85").Value = 8513.25
$ws.Range("M85").Value = -7265.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2237.125
$ws.Range("I113").Value = 2237.125
$ws.Range("K113").Value = 2237.125
$ws.Range("M113").Value = -67.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 32498.125
$ws.Range("I122").Value = 32498.125
$ws.Range("K122").Value = 97494.375
$ws.Range("M122").Value = -95044.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 7716.04
$ws.Range("I136").Value = 7647.15
$ws.Range("J136").Value = 7991.6
$ws.Range("K136").Value = 22941.45
$ws.Range("L136").Value = 23974.8
$ws.Range("M136").Value = -20391.45
$ws.Range("N136").Value = -29074.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 758.381
$ws.Range("I107").Value = 853.75
$ws.Range("J107").Value = 699.6923
$ws.Range("K107").Value = 2561.25
$ws.Range("L107").Value = 2099.0769
$ws.Range("M107").Value = -641.25
$ws.Range("N107").Value = -5939.0769

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1566.5264
$ws.Range("I113").Value = 941.5625
$ws.Range("K113").Value = 2824.6875
$ws.Range("M113").Value = -654.6875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 11156.682
$ws.Range("I136").Value = 12011.871
$ws.Range("K136").Value = 36035.613
$ws.Range("M136").Value = -33485.613
